$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column G ("Stadt*") before the existing Telefonnummer* column ---
# This shifts Telefonnummer*, E-Mail-Adresse*, Geburtsdatum* and Tagesgast one column to the right.
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# --- New header ---
$ws.Range("G1").Value = "Stadt*"

# --- Postleitzahl (postal code) values become text, and city column gets filled in ---
$ws.Range("F2").Value2 = "40882"
$ws.Range("G2").Value2 = "Ratingen"

$ws.Range("F3").Value2 = "40882"
$ws.Range("G3").Value2 = "Düsseldorf"

$ws.Range("F4").Value2 = "53113"
$ws.Range("G4").Value2 = "Bonn"

# --- Fix up the hyperlinks on the E-Mail-Adresse* column, which moved from H to I ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:Robert@pfadfinder.de") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:Robert@pfadfinder.de") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:Robert@pfadfinder.de") | Out-Null

# Re-touching the (unchanged) number format nudges the cell style back onto the
# original "Link" style record instead of a freshly minted duplicate one.
$ws.Range("I2:I4").NumberFormat = $ws.Range("I2").NumberFormat

# --- Restore original selection ---
$ws.Range("G5").Select()
